# Add "Insertion sort" row to the Problem-soln sheet, mirroring the
# existing Selection sort / Bubble sort rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Problem-soln")

$insertionUrl = "https://www.codingninjas.com/studio/problems/insertion-sort_624381?utm_source=striver&utm_medium=website&utm_campaign=a_zcoursetuf&leftPanelTab=1"
$insertionDesc = "keep comparing every element with left elements , keep swapping until left is greater`n"

# New row 11: Insertion sort / link / description (set before C10 so the
# shared-string table order matches: Insertion sort, link, description,
# then "adjcent swapping ").
$ws.Range("A11").Value = "Insertion sort"
$ws.Range("B11").Value = $insertionUrl
$ws.Range("C11").Value = $insertionDesc

$ws.Hyperlinks.Add($ws.Range("B11"), $insertionUrl) | Out-Null

$ws.Range("A11").Style = $ws.Range("A10").Style
$ws.Range("B11").Style = $ws.Range("B10").Style
$ws.Range("C11").Style = $ws.Range("A10").Style
$ws.Rows.Item(11).RowHeight = $ws.Rows.Item(10).RowHeight

# Row 10, column C currently empty -> "adjcent swapping "
$ws.Range("C10").Value = "adjcent swapping "
$ws.Range("C10").Style = $ws.Range("A10").Style

$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("F10").Select()
